# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.401.96"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.59"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.41"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3759"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.87"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3423"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.004"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.30"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.31"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06727"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.73"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.229"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.403.18"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.648"
$ws.Range("E26").Value = "  -10.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.16"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.96"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.022"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.68"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.745.72"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.158"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9821"
$ws.Range("E34").Value = "  -5.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.11"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08485"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02538"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.371"
$ws.Range("E38").Value = "  +10.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2315"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06546"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.418"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.45"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6379"
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.08"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5973"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.291"
$ws.Range("E48").Value = "  +3.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.094"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.84"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07332"
$ws.Range("E51").Value = "  +0.99%  "
